$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Add a header row above the existing Edwards/Govindaraju compound-weight
# data so the two columns are clearly labeled.
$ws.Range("A1").Value = "Compound"
$ws.Range("B1").Value = "Weight"

# Restore the view: scrolled down with D15 as the active selection.
$ws.Range("D15").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
